$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing contents (keeps formatting/styles) so the shared-string table
# is rebuilt cleanly in the exact order the new values are written below.
$ws.Cells.ClearContents()

$ws.Range("A1").Value = "Feature"
$ws.Range("B1").Value = "Importance"

$ws.Range("A2").Value = 'sob código'
$ws.Range("B2").Value = 0.0625
$ws.Range("A3").Value = 'concessão ordem'
$ws.Range("B3").Value = 0.0625
$ws.Range("A4").Value = 'xliii constituição'
$ws.Range("B4").Value = 0.03125
$ws.Range("A5").Value = 'relatório http'
$ws.Range("B5").Value = 0.03125
$ws.Range("A6").Value = 'deferi'
$ws.Range("B6").Value = 0.03125
$ws.Range("A7").Value = 'deferida'
$ws.Range("B7").Value = 0.03125
$ws.Range("A8").Value = 'deferido'
$ws.Range("B8").Value = 0.03125
$ws.Range("A9").Value = 'deferimento'
$ws.Range("B9").Value = 0.03125
$ws.Range("A10").Value = 'efeitos ordem'
$ws.Range("B10").Value = 0.03125
$ws.Range("A11").Value = 'outro motivo'
$ws.Range("B11").Value = 0.03125
$ws.Range("A12").Value = 'ordem prisão'
$ws.Range("B12").Value = 0.03125
$ws.Range("A13").Value = 'opina deferimento'
$ws.Range("B13").Value = 0.03125
$ws.Range("A14").Value = 'opina concessão'
$ws.Range("B14").Value = 0.03125
$ws.Range("A15").Value = 'ministro gilson'
$ws.Range("B15").Value = 0.03125
$ws.Range("A16").Value = 'medidas cautelares'
$ws.Range("B16").Value = 0.03125
$ws.Range("A17").Value = 'liminar suspender'
$ws.Range("B17").Value = 0.03125
$ws.Range("A18").Value = 'liminar hc'
$ws.Range("B18").Value = 0.03125
$ws.Range("A19").Value = 'gabinete prestou'
$ws.Range("B19").Value = 0.03125
$ws.Range("A20").Value = 'justiça indeferiu'
$ws.Range("B20").Value = 0.03125
$ws.Range("A21").Value = 'inconstitucional'
$ws.Range("B21").Value = 0.03125
$ws.Range("A22").Value = 'república concessão'
$ws.Range("B22").Value = 0.03125
$ws.Range("A23").Value = 'indeferiu liminarmente'
$ws.Range("B23").Value = 0.03125
$ws.Range("A24").Value = 'senha relatório'
$ws.Range("B24").Value = 0.03125
$ws.Range("A25").Value = 'concessão'
$ws.Range("B25").Value = 0.03125
$ws.Range("A26").Value = 'cento'
$ws.Range("B26").Value = 0.03125
$ws.Range("A27").Value = 'violência grave'
$ws.Range("B27").Value = 0.03125
$ws.Range("A28").Value = 'vedação liberdade'
$ws.Range("B28").Value = 0.03125
$ws.Range("A29").Value = 'restritiva direitos'
$ws.Range("B29").Value = 0.03125
$ws.Range("A30").Value = 'submetida'
$ws.Range("B30").Value = 0.03125
$ws.Range("A31").Value = 'sobrestamento'
$ws.Range("B31").Value = 0.03125
$ws.Range("A32").Value = 'suspender'
$ws.Range("B32").Value = 0
$ws.Range("A33").Value = 'mérito deste'
$ws.Range("B33").Value = 0
$ws.Range("A34").Value = 'substituição pena'
$ws.Range("B34").Value = 0
$ws.Range("A35").Value = 'superação súmula'
$ws.Range("B35").Value = 0
$ws.Range("A36").Value = 'liminar deferida'
$ws.Range("B36").Value = 0
$ws.Range("A37").Value = 'suspender efeitos'
$ws.Range("B37").Value = 0
$ws.Range("A38").Value = 'liminar espécie'
$ws.Range("B38").Value = 0
$ws.Range("A39").Value = 'liminar assessor'
$ws.Range("B39").Value = 0
$ws.Range("A40").Value = 'liberdade restritiva'
$ws.Range("B40").Value = 0
$ws.Range("A41").Value = 'stj indeferiu'
$ws.Range("B41").Value = 0
$ws.Range("A42").Value = 'restritiva'
$ws.Range("B42").Value = 0
$ws.Range("A43").Value = 'opinou concessão'
$ws.Range("B43").Value = 0
$ws.Range("A44").Value = 'requer medida'
$ws.Range("B44").Value = 0
$ws.Range("A45").Value = 'senha primeira'
$ws.Range("B45").Value = 0
$ws.Range("A46").Value = 'precário'
$ws.Range("B46").Value = 0
$ws.Range("A47").Value = 'senha'
$ws.Range("B47").Value = 0
$ws.Range("A48").Value = 'precário efêmero'
$ws.Range("B48").Value = 0
$ws.Range("A49").Value = 'revelou contornos'
$ws.Range("B49").Value = 0
$ws.Range("A50").Value = 'resumida prisão'
$ws.Range("B50").Value = 0
$ws.Range("A51").Value = 'preventiva fundamentos'
$ws.Range("B51").Value = 0
$ws.Range("A52").Value = 'resumida'
$ws.Range("B52").Value = 0
$ws.Range("A53").Value = 'processo formalizado'
$ws.Range("B53").Value = 0
$ws.Range("A54").Value = 'informado análise'
$ws.Range("B54").Value = 0
$ws.Range("A55").Value = 'pública estado'
$ws.Range("B55").Value = 0
$ws.Range("A56").Value = 'proferida ministro'
$ws.Range("B56").Value = 0
$ws.Range("A57").Value = 'assuntos FIANÇA'
$ws.Range("B57").Value = 0
$ws.Range("A58").Value = 'informado'
$ws.Range("B58").Value = 0
$ws.Range("A59").Value = 'inadmissão'
$ws.Range("B59").Value = 0
$ws.Range("A60").Value = 'contornos impetração'
$ws.Range("B60").Value = 0
$ws.Range("A61").Value = 'contornos'
$ws.Range("B61").Value = 0
$ws.Range("A62").Value = 'cautelares previstas'
$ws.Range("B62").Value = 0
$ws.Range("A63").Value = 'campo precário'
$ws.Range("B63").Value = 0
$ws.Range("A64").Value = 'aurélio decisão'
$ws.Range("B64").Value = 0
$ws.Range("A65").Value = 'assim revelou'
$ws.Range("B65").Value = 0
$ws.Range("A66").Value = 'assim resumida'
$ws.Range("B66").Value = 0
$ws.Range("A67").Value = 'arquivado'
$ws.Range("B67").Value = 0
$ws.Range("A68").Value = 'análise pedido'
$ws.Range("B68").Value = 0
$ws.Range("A69").Value = 'alcance'
$ws.Range("B69").Value = 0
$ws.Range("A70").Value = 'afastamento enunciado'
$ws.Range("B70").Value = 0
$ws.Range("A71").Value = 'Relator_OCTAVIO GALLOTTI'
$ws.Range("B71").Value = 0
$ws.Range("A72").Value = 'Relator_ILMAR GALVÃO'
$ws.Range("B72").Value = 0
$ws.Range("A73").Value = 'assuntos TRANSFERÊNCIA DE PRESO'
$ws.Range("B73").Value = 0
$ws.Range("A74").Value = 'assuntos FURTO (ART. 155)'
$ws.Range("B74").Value = 0
$ws.Range("A75").Value = 'assuntos LICITAÇÕES'
$ws.Range("B75").Value = 0
$ws.Range("A76").Value = 'assuntos DESCLASSIFICAÇÃO'
$ws.Range("B76").Value = 0
$ws.Range("A77").Value = 'assuntos INDEFERIMENTO'
$ws.Range("B77").Value = 0
$ws.Range("A78").Value = 'assuntos COMPETÊNCIA DO MP'
$ws.Range("B78").Value = 0
$ws.Range("A79").Value = 'código senha'
$ws.Range("B79").Value = 0
$ws.Range("A80").Value = 'decisão proferida'
$ws.Range("B80").Value = 0
$ws.Range("A81").Value = 'decretada desfavor'
$ws.Range("B81").Value = 0
$ws.Range("A82").Value = 'ficou'
$ws.Range("B82").Value = 0
$ws.Range("A83").Value = 'implicou deferimento'
$ws.Range("B83").Value = 0
$ws.Range("A84").Value = 'assuntos HABEAS CORPUS - LIBERATÓRIO'
$ws.Range("B84").Value = 0
$ws.Range("A85").Value = 'http sob'
$ws.Range("B85").Value = 0
$ws.Range("A86").Value = 'http'
$ws.Range("B86").Value = 0
$ws.Range("A87").Value = 'fundamentos insubsistência'
$ws.Range("B87").Value = 0
$ws.Range("A88").Value = 'formalizado ato'
$ws.Range("B88").Value = 0
$ws.Range("A89").Value = 'flagrante preventiva'
$ws.Range("B89").Value = 0
$ws.Range("A90").Value = 'ficou assim'
$ws.Range("B90").Value = 0
$ws.Range("A91").Value = 'estariam ausentes'
$ws.Range("B91").Value = 0
$ws.Range("A92").Value = 'deduzida'
$ws.Range("B92").Value = 0
$ws.Range("A93").Value = 'espécie ficou'
$ws.Range("B93").Value = 0
$ws.Range("A94").Value = 'enunciado súmula'
$ws.Range("B94").Value = 0
$ws.Range("A95").Value = 'enunciado'
$ws.Range("B95").Value = 0
$ws.Range("A96").Value = 'eis informado'
$ws.Range("B96").Value = 0
$ws.Range("A97").Value = 'efêmero'
$ws.Range("B97").Value = 0
$ws.Range("A98").Value = 'deserção'
$ws.Range("B98").Value = 0
$ws.Range("A99").Value = 'deferida assessoria'
$ws.Range("B99").Value = 0
$ws.Range("A100").Value = 'deferi pedido'
$ws.Range("B100").Value = 0
$ws.Range("A101").Value = 'impetração eis'
$ws.Range("B101").Value = 0

Write-Host "Feature importance sheet updated"
